# Trade #9 closed at 2026-02-16 21:53:09 - leadlag UP +0.000%
# Append the new (still-OPEN) trade row to both the "All Trades" sheet and
# the strategy-specific "leadlag" sheet.

$wb = $excel.ActiveWorkbook

function Add-TradeRow($ws, $row) {
    # Trade #
    $ws.Cells.Item($row, 1).Value = 9

    # Date / Time must stay as literal text, not get auto-converted to an
    # Excel date/time serial number. Forcing the cell to "Text" format
    # before assignment prevents the automatic date parsing; resetting the
    # style back to "Normal" afterwards keeps the cell style clean (same
    # default style as the rest of the sheet).
    $cB = $ws.Cells.Item($row, 2)
    $cB.NumberFormat = "@"
    $cB.Value = "2026-02-16"
    $cB.Style = "Normal"

    $cC = $ws.Cells.Item($row, 3)
    $cC.NumberFormat = "@"
    $cC.Value = "21:53:09"
    $cC.Style = "Normal"

    $ws.Cells.Item($row, 4).Value = "leadlag"
    $ws.Cells.Item($row, 5).Value = "UP"
    $ws.Cells.Item($row, 6).Value = 68293.08

    # Exit Price: trade is still OPEN, so this stays blank. Touch the cell
    # (same text-format trick as above) so a placeholder cell element is
    # still written out for this column, same as the other blank cells.
    $cG = $ws.Cells.Item($row, 7)
    $cG.NumberFormat = "@"
    $cG.Value = ""
    $cG.Style = "Normal"

    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0.617
    $ws.Cells.Item($row, 13).Value = "Binance leading with 0.062% move"

    # Exit Reason: also blank while the trade remains OPEN.
    $cN = $ws.Cells.Item($row, 14)
    $cN.NumberFormat = "@"
    $cN.Value = ""
    $cN.Style = "Normal"

    $ws.Cells.Item($row, 15).Value = 0
}

$wsAll = $wb.Worksheets.Item("All Trades")
Add-TradeRow $wsAll 10

$wsLeadlag = $wb.Worksheets.Item("leadlag")
Add-TradeRow $wsLeadlag 9
